# Rebuild the director/movie/gross table using the full dataset (found via
# looking up each director's children/filmography) instead of the previous
# one-movie-per-director placeholder rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same.
$ws.Cells.Item(1, 1).Value = "Director"
$ws.Cells.Item(1, 2).Value = "Movie"
$ws.Cells.Item(1, 3).Value = "Worldwide Gross"

# Director, Movie, Worldwide Gross - sorted by gross descending.
$movieData = @(
    @("james cameron", "Avatar", 2920357254),
    @("james cameron", "Titanic", 2201647264),
    @("steven speilberg", "Jurassic Park", 1109802321),
    @("christopher nolan", "The Dark Knight Rises", 1081153097),
    @("christopher nolan", "The Dark Knight", 1006234167),
    @("taika waititi", "Thor: Ragnarok", 853983879),
    @("christopher nolan", "Inception", 836848102),
    @("steven speilberg", "E.T. the Extra-Terrestrial", 792910554),
    @("steven speilberg", "Indiana Jones and the Kingdom of the Crystal Skull", 790653942),
    @("taika waititi", "Thor: Love and Thunder", 760677374),
    @("christopher nolan", "Interstellar", 716218351),
    @("steven speilberg", "The Lost World: Jurassic Park", 618638999),
    @("steven speilberg", "War of the Worlds", 603873119),
    @("christopher nolan", "Dunkirk", 527016307),
    @("james cameron", "Terminator 2: Judgment Day", 520881154),
    @("james cameron", "True Lies", 378882411),
    @("james cameron", "Aliens", 131060248),
    @("taika waititi", "Jojo Rabbit", 90335025),
    @("taika waititi", "Boy", 43551154),
    @("taika waititi", "Hunt for the Wilderpeople", 23915910)
)

$row = 2
foreach ($movie in $movieData) {
    $ws.Cells.Item($row, 1).Value = $movie[0]
    $ws.Cells.Item($row, 2).Value = $movie[1]
    $ws.Cells.Item($row, 3).Value = $movie[2]
    $row = $row + 1
}
